$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 414
$ws.Range("F5").Value = 51
$ws.Range("F8").Value = 1077
$ws.Range("F10").Value = 374
$ws.Range("F11").Value = 431
$ws.Range("F18").Value = 561
$ws.Range("F20").Value = 5725
$ws.Range("F22").Value = 1612
$ws.Range("F24").Value = 60
$ws.Range("F26").Value = 5312
$ws.Range("F27").Value = 5312
$ws.Range("F33").Value = 57
$ws.Range("F34").Value = 44
$ws.Range("F36").Value = 106
$ws.Range("F38").Value = 3810

# Sheet "演出" (Performances)
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F4").Value = 34
$ws.Range("F5").Value = 160
$ws.Range("F8").Value = 177

# Sheet "本地生活" (Local Life)
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 9416
$ws.Range("F4").Value = 2155

# Sheet "全部类型" (All Types)
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 9416
$ws.Range("F4").Value = 2155
$ws.Range("F6").Value = 414
$ws.Range("F8").Value = 51
$ws.Range("F11").Value = 1077
$ws.Range("F12").Value = 374
$ws.Range("F13").Value = 431
$ws.Range("F22").Value = 5725
$ws.Range("F24").Value = 1612
$ws.Range("F30").Value = 5312
$ws.Range("F31").Value = 5312
$ws.Range("F37").Value = 44
$ws.Range("F39").Value = 106
$ws.Range("F47").Value = 3810

$wb.Save()
